$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: "LearnPortal: " label + hyperlinked tutorial title
$ws.Range("B11").Value = "LearnPortal: "

$ws.Hyperlinks.Add(
    $ws.Range("D11"),
    "https://randomnerdtutorials.com/esp32-mpu-6050-accelerometer-gyroscope-arduino/",
    "",
    "",
    "https://randomnerdtutorials.com/esp32-mpu-6050-accelerometer-gyroscope-arduino/"
)

# Restore the friendly display text + hyperlink styling for the cell
$ws.Range("D11").Value = "ESP32 MPU-6050 Accelerometer and Gyroscope (Arduino) | Random Nerd Tutorials"
$ws.Range("D11").Style = "Link"

$ws.Range("D11").Select()
